$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")

# B11 currently holds the text "R40" (a shared string). The target value is
# the text "1" -- which LOOKS numeric, so a plain .Value assignment would be
# auto-coerced by Excel into a real number (and changing NumberFormat to
# force text on the cell itself would also permanently change its style
# index). Instead, build the text "1" on a scratch cell via a formula (so it
# is unambiguously a string), then paste only the VALUE into B11 - leaving
# B11's existing style/format completely untouched.
$helper = $ws.Range("Z1")
$helper.Formula = '="1"'

$helper.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues

$helper.Clear() | Out-Null
$excel.CutCopyMode = 0

$wb.Save()
